$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("FPIEBP")

# About sheet: update the "last updated" date in C1 (1/3/24 -> 3/28/24)
$ws1.Range("C1").Value = 45379

# FPIEBP sheet: re-rank the priorities for "hard coal" (row 3)
$ws2.Range("B3").Value = 1
$ws2.Range("C3").Value = 3
$ws2.Range("D3").Value = 2

# Leave FPIEBP as the active sheet, with its selection moved to E3
$ws2.Activate()
$ws2.Range("E3").Select()
